$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7005266547203064
$ws.Range("B1").Value = 1.33527147769928
$ws.Range("C1").Value = 3.778372526168823
$ws.Range("D1").Value = 2.569887161254883
$ws.Range("E1").Value = 0.6293447613716125
